$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 holds the "last updated" timestamp banner.
$ws.Range("A1").Value = "Datos actualizados a 10 de Septiembre de 2020 a las 17:25"

$csv = @"
Estados Unidos|6556486|7011|3856335|2504623|0|289|195528
India|4494389|31424|3490908|928153|0|237|75328
Brasil|4199332|0|3453336|617343|0|0|128653
Rusia|1046370|5363|862373|165734|0|128|18263
Peru|702776|0|536959|135581|0|0|30236
Colombia|686856|0|552885|111918|0|0|22053
Mexico|647507|4647|454982|123430|0|611|69095
Sudafrica|642431|0|569935|57328|0|0|15168
España|543379|0|0|0|0|0|29628
Argentina|512293|0|390098|111482|0|55|10713
Chile|428669|1642|401356|15532|0|79|11781
Iran|395488|2063|340842|31848|0|129|22798
Reino Unido|355219|0|0|0|0|0|41594
Francia|344101|0|88524|224783|0|0|30794
Banglades|332970|1892|233550|94786|0|41|4634
Arabia Saudita|323720|708|299998|19533|0|24|4189
Pakistan|299855|196|287950|5540|0|6|6365
Turquia|284943|0|254188|23918|0|0|6837
Italia|283180|1597|211885|35708|0|10|35587
Irak|278418|4597|213817|56787|0|82|7814
Alemania|256835|486|231900|15523|0|2|9412
Filipinas|248947|3821|186058|58823|0|80|4066
Indonesia|207203|3861|147510|51237|0|120|8456
Ucrania|145612|2582|65877|76712|0|44|3023
Israel|143049|1952|109775|32219|0|1|1055
Canada|134464|170|118271|7037|0|1|9156
Bolivia|123345|1037|77216|38983|0|49|7146
Catar|121052|206|117978|2869|0|0|205
Ecuador|112166|0|91242|10223|0|0|10701
Kazajistan|106584|86|100248|4702|0|0|1634
Republica Dominicana|100937|0|74305|24718|0|0|1914
Egipto|100403|0|80689|14137|0|0|5577
Rumania|99684|1380|41010|54609|0|47|4065
Panama|99042|0|71419|25507|0|0|2116
Kuwait|92822|740|83040|9226|0|4|556
Belgica|89691|550|18635|61139|0|5|9917
Oman|88337|398|83325|4250|0|11|762
Suecia|86194|0|0|0|0|4|5843
China|85153|7|80358|161|0|0|4634
Guatemala|79622|0|68308|8417|0|0|2897
Paises Bajos|78511|823|0|0|0|3|6249
Marruecos|77878|0|59723|16702|0|0|1453
Emiratos Arabes Unidos|76911|930|67945|8568|0|5|398
Bielorrusia|73591|189|72203|656|0|6|732
Japon|72726|0|64100|7233|0|0|1393
Polonia|72453|506|58069|12225|0|12|2159
Honduras|65597|379|15064|48489|0|10|2044
Portugal|62126|585|43441|16833|0|3|1852
Etiopia|61700|0|23054|37680|0|0|966
Barein|57450|0|51819|5425|0|3|206
Singapur|57229|63|56492|710|0|0|27
Venezuela|56751|0|45318|10981|0|0|452
Nigeria|55632|0|43610|10952|0|0|1070
Costa Rica|51224|0|19544|31137|0|0|543
Nepal|50465|1246|35700|14448|0|5|317
Argelia|47216|0|33379|12256|0|0|1581
Suiza|45711|405|38100|5591|0|1|2020
Uzbekistan|45473|543|42509|2594|0|4|370
Armenia|45326|174|41233|3187|0|1|906
Ghana|45313|0|44188|842|0|0|283
Kirguistan|44684|71|40487|3136|0|0|1061
Moldavia|41704|560|28962|11636|0|10|1106
Afganistan|38572|28|31129|6023|0|0|1420
Azerbaiyan|37874|142|35269|2049|0|1|556
Kenia|35460|0|21557|13296|0|0|607
Serbia|32136|58|30943|464|0|1|729
Austria|31247|664|26043|4456|0|1|748
Chequia|31036|0|20374|10218|0|0|444
Irlanda|30164|0|23364|5019|0|0|1781
Estado de Palestina|28664|745|18821|9645|0|6|198
El Salvador|26688|86|17106|8808|0|4|774
Australia|26513|48|23036|2689|0|7|788
Paraguay|25026|0|12309|12243|0|0|474
Bosnia y Herzegovina|22544|286|15273|6591|0|5|680
Libano|21877|0|7024|14641|0|0|212
Corea del Sur|21743|155|17360|4037|0|2|346
Libia|20939|477|2420|18180|0|15|339
Camerun|19848|0|18448|985|0|0|415
Dinamarca|18924|317|16069|2226|0|1|629
Costa de Marfil|18815|0|17770|926|0|0|119
Bulgaria|17435|0|12474|4259|0|0|702
Madagascar|15624|104|14295|1121|0|2|208
Republica de Macedonia|15293|0|12754|1905|0|0|634
Senegal|14150|48|10207|3650|0|0|293
Sudan|13437|0|6730|5874|0|0|833
Zambia|13214|102|11899|1015|0|0|300
Croacia|12917|291|10142|2567|0|2|208
Grecia|12080|0|3804|7983|0|0|293
Noruega|11801|55|9348|2188|0|1|265
Albania|10704|0|6284|4098|0|0|322
Consejo Danes para los Refugiados|10343|19|9533|548|0|2|262
Hungria|10191|476|3990|5571|0|2|630
Guinea|9885|0|9068|754|0|0|63
Malasia|9628|45|9167|333|0|0|128
Guayana Francesa|9418|0|8986|370|0|0|62
Namibia|9256|148|4767|4393|0|3|96
Tayikistan|8939|40|7714|1153|0|1|72
Maldivas|8834|0|6288|2515|0|0|31
Gabon|8621|0|7618|950|0|0|53
Finlandia|8469|39|7500|632|0|0|337
Haiti|8384|0|5991|2179|0|0|214
Zimbabue|7429|0|5542|1665|0|0|222
Mauritania|7191|0|6701|329|0|0|161
Luxemburgo|7023|0|6292|607|0|0|124
Montenegro|6094|0|4350|1630|0|0|114
Malaui|5653|0|3630|1847|0|0|176
Tunez|5417|0|1862|3459|0|0|96
Republica de Yibuti|5391|0|5327|3|0|0|61
Eslovaquia|5066|178|3001|2028|0|0|37
Guinea Ecuatorial|4990|0|4479|428|0|0|83
Suazilandia|4936|0|4077|762|0|0|97
Hong Kong|4914|12|4582|233|0|0|99
Congo|4891|0|3887|921|0|0|83
Nicaragua|4818|0|2913|1761|0|0|144
Mozambique|4764|0|2763|1973|0|0|28
Republica de Africa Central|4736|0|1825|2849|0|0|62
Cabo Verde|4473|0|3915|515|0|0|43
Ruanda|4460|0|2325|2114|0|0|21
Cuba|4459|0|3727|628|0|0|104
Surinam|4447|0|3622|733|0|0|92
Uganda|4291|190|1998|2245|0|2|48
Tailandia|3454|7|3310|86|0|0|58
Jamaica|3437|114|1019|2380|0|2|38
Eslovenia|3389|77|2620|634|0|0|135
Mayotte|3374|0|2964|370|0|0|40
Somalia|3371|0|2738|536|0|0|97
Siria|3351|0|780|2428|0|0|143
Gambia|3293|0|1460|1734|0|0|99
Lituania|3199|36|2030|1083|0|0|86
Sri Lanka|3152|5|2955|185|0|0|12
Angola|3092|0|1245|1721|0|0|126
Mali|2898|0|2267|503|0|0|128
Jordania|2739|80|1981|738|0|1|20
Aruba|2730|0|1351|1364|0|0|15
Bahamas|2721|64|1183|1475|0|0|63
Estonia|2600|15|2223|313|0|0|64
Trinidad yTobago|2588|0|749|1800|0|0|39
Sudan del Sur|2555|0|1290|1216|0|0|49
Reunion|2416|0|1313|1090|0|0|13
Guadalupe|2287|0|336|1928|0|0|23
Guinea-Bisau|2245|0|1127|1080|0|0|38
Benin|2242|0|1793|409|0|0|40
Malta|2204|42|1803|387|0|0|14
Islandia|2153|0|2067|76|0|0|10
Birmania|2150|261|625|1511|0|2|14
Botsuana|2126|0|493|1624|0|0|9
Sierra Leona|2067|0|1622|373|0|0|72
Yemen|1999|0|1209|214|0|0|576
Georgia|1830|57|1334|477|0|0|19
Nueva Zelanda|1792|4|1648|120|0|0|24
Uruguay|1741|0|1478|218|0|0|45
Guyana|1703|0|1045|610|0|0|48
Togo|1528|0|1144|348|0|0|36
Republica de Chipre|1514|0|1237|255|0|0|22
Burkina Faso|1476|0|1118|302|0|0|56
Letonia|1448|5|1248|165|0|0|35
Belice|1365|4|326|1021|0|2|18
Liberia|1313|0|1194|37|0|0|82
Principado de Andorra|1301|0|938|310|0|0|53
Niger|1178|0|1099|10|0|0|69
Lesoto|1164|0|528|605|0|0|31
Vietnam|1059|0|890|134|0|0|35
Republica del Chad|1048|0|931|38|0|0|79
Martinica|939|0|98|823|0|0|18
Santo Tome y Principe|898|0|862|21|0|0|15
Polinesia Francesa|795|0|528|267|0|0|0
San Marino|722|0|662|18|0|0|42
Crucero|712|0|651|48|0|0|13
Islas Turcas y Caicos|628|14|266|357|0|0|5
San Martin (Parte Holandesa)|530|3|333|178|0|0|19
Tanzania|509|0|183|305|0|0|21
Papua Nueva Guinea|504|1|232|267|0|0|5
Taiwan|496|1|475|14|0|0|7
Burundi|466|0|374|91|0|0|1
Comoras|456|0|415|34|0|0|7
Islas Feroe|415|0|410|5|0|0|0
Mauricio|361|0|335|16|0|0|10
Eritrea|341|0|295|46|0|0|0
Isla de Man|337|0|312|1|0|0|24
Gibraltar|323|1|287|36|0|0|0
Mongolia|310|0|298|12|0|0|0
Camboya|274|0|273|1|0|0|0
San Martin (Parte Francesa)|256|0|107|143|0|0|6
Butan|234|0|155|79|0|0|0
Islas Caimanes|207|0|204|2|0|0|1
Barbados|180|0|156|17|0|0|7
Bermudas|177|0|159|9|0|0|9
Monaco|161|0|107|53|0|0|1
Brunei|145|0|139|3|0|0|3
Seychelles|137|0|127|10|0|0|0
Curazao|116|0|51|64|0|0|1
Liechtenstein|109|1|105|3|0|0|1
Antigua y Barbuda|95|0|91|1|0|0|3
Islas Virgenes Britanicas|63|0|10|52|0|0|1
San Vicente y las Granadinas|62|0|58|4|0|0|0
Macao|46|0|46|0|0|0|0
Puerto Rico|39|0|1|36|0|0|2
Guam|32|0|0|31|0|0|1
Fiyi|32|0|24|6|0|0|2
Timor Oriental|27|0|25|2|0|0|0
Santa Lucia|27|0|26|1|0|0|0
Nueva Caledonia|26|0|25|1|0|0|0
Bonaire, San Eustaquio y Saba|25|4|7|18|0|0|0
Dominica|24|2|18|6|0|0|0
Granada|24|0|24|0|0|0|0
Laos|22|0|21|1|0|0|0
San Bartolome|21|0|13|8|0|0|0
Islas Virgenes de los Estados Unidos|17|0|0|17|0|0|0
San Cristobal y Nieves|17|0|17|0|0|0|0
Groenlandia|14|0|14|0|0|0|0
Islas Malvinas|13|0|13|0|0|0|0
Montserrat|13|0|12|0|0|0|1
Santa Sede|12|0|12|0|0|0|0
San Pedro y Miquelon|10|1|5|5|0|0|0
Sahara Occidental|10|0|8|1|0|0|1
Anguila|3|0|3|0|0|0|0
"@

$lines = $csv -split "`r?`n" | Where-Object { $_.Length -gt 0 }
$rowCount = $lines.Count
$colCount = 8

$data = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $fields = $lines[$i] -split '\|'
    $data[$i, 0] = $fields[0]
    for ($j = 1; $j -lt $colCount; $j++) {
        $data[$i, $j] = [long]$fields[$j]
    }
}

$startRow = 4
$endRow = $startRow + $rowCount - 1
$range = $ws.Range("A" + $startRow + ":H" + $endRow)
$range.Value = $data
